$d = $word.ActiveDocument

$d.Content.Find.Execute("730÷6=121, 4", $true, $true, $false, $false, $false, $true, 1, $false, "785÷6=130, 5", 2)
$d.Content.Find.Execute("890÷6=148, 2", $true, $true, $false, $false, $false, $true, 1, $false, "387÷9=43, 0", 2)
$d.Content.Find.Execute("690÷6=115, 0", $true, $true, $false, $false, $false, $true, 1, $false, "397÷2=198, 1", 2)
$d.Content.Find.Execute("270÷3=90, 0", $true, $true, $false, $false, $false, $true, 1, $false, "143÷3=47, 2", 2)
$d.Content.Find.Execute("513÷9=57, 0", $true, $true, $false, $false, $false, $true, 1, $false, "418÷4=104, 2", 2)
$d.Content.Find.Execute("133÷9=14, 7", $true, $true, $false, $false, $false, $true, 1, $false, "484÷8=60, 4", 2)
$d.Content.Find.Execute("693÷9=77, 0", $true, $true, $false, $false, $false, $true, 1, $false, "956÷5=191, 1", 2)
$d.Content.Find.Execute("495÷5=99, 0", $true, $true, $false, $false, $false, $true, 1, $false, "660÷5=132, 0", 2)
$d.Content.Find.Execute("572÷2=286, 0", $true, $true, $false, $false, $false, $true, 1, $false, "675÷6=112, 3", 2)
$d.Content.Find.Execute("149÷6=24, 5", $true, $true, $false, $false, $false, $true, 1, $false, "396÷3=132, 0", 2)
$d.Content.Find.Execute("132÷8=16, 4", $true, $true, $false, $false, $false, $true, 1, $false, "881÷7=125, 6", 2)
$d.Content.Find.Execute("985÷3=328, 1", $true, $true, $false, $false, $false, $true, 1, $false, "841÷6=140, 1", 2)
$d.Content.Find.Execute("389÷7=55, 4", $true, $true, $false, $false, $false, $true, 1, $false, "424÷9=47, 1", 2)
$d.Content.Find.Execute("169÷8=21, 1", $true, $true, $false, $false, $false, $true, 1, $false, "401÷5=80, 1", 2)
$d.Content.Find.Execute("755÷2=377, 1", $true, $true, $false, $false, $false, $true, 1, $false, "877÷5=175, 2", 2)
$d.Content.Find.Execute("198÷4=49, 2", $true, $true, $false, $false, $false, $true, 1, $false, "807÷8=100, 7", 2)
$d.Content.Find.Execute("155÷2=77, 1", $true, $true, $false, $false, $false, $true, 1, $false, "823÷8=102, 7", 2)
$d.Content.Find.Execute("959÷3=319, 2", $true, $true, $false, $false, $false, $true, 1, $false, "984÷2=492, 0", 2)
$d.Content.Find.Execute("182÷5=36, 2", $true, $true, $false, $false, $false, $true, 1, $false, "410÷4=102, 2", 2)
$d.Content.Find.Execute("900÷7=128, 4", $true, $true, $false, $false, $false, $true, 1, $false, "249÷3=83, 0", 2)
$d.Content.Find.Execute("321÷2=160, 1", $true, $true, $false, $false, $false, $true, 1, $false, "799÷7=114, 1", 2)
$d.Content.Find.Execute("178÷6=29, 4", $true, $true, $false, $false, $false, $true, 1, $false, "407÷7=58, 1", 2)
$d.Content.Find.Execute("609÷5=121, 4", $true, $true, $false, $false, $false, $true, 1, $false, "812÷9=90, 2", 2)
$d.Content.Find.Execute("923÷5=184, 3", $true, $true, $false, $false, $false, $true, 1, $false, "624÷3=208, 0", 2)
$d.Content.Find.Execute("683÷8=85, 3", $true, $true, $false, $false, $false, $true, 1, $false, "174÷6=29, 0", 2)
